$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "64.540.53"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.41%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.159.19"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.63%  "
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "614.30"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.12%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.00"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.91%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.156.88"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.60%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.525"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.35%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.151"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.08%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.41"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.54%  "
$ws.Range("E12").Value = "  -0.87%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000257"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.22%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.57"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.24%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.678.99"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.81%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "64.563.73"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.47%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.158.88"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.70%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.86"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.02%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "478.80"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.11%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.60"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.37%  "
$ws.Range("E22").Value = "  +2.33%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.94"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.92%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.80"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.69%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "83.81"
$ws.Range("D25").Style = "Normal"
$ws.Range("E26").Value = "  +0.00%  "
$ws.Range("E27").Value = "  -3.27%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.62"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.63%  "
$ws.Range("E29").Value = "  +3.80%  "
$ws.Range("E30").Value = "  -3.86%  "
$ws.Range("E31").Value = "  -5.31%  "
$ws.Range("E32").Value = "  +0.11%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "26.46"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.32%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.66"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.80%  "
$ws.Range("E35").Value = "  +2.10%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0₃0778"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +7.46%  "
$ws.Range("E37").Value = "  -0.92%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "53.15"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.32%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.15"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.30%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "461.48"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.14%  "
$ws.Range("E42").Value = "  -3.42%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.32"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.07%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.856.20"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.60%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.30"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.16%  "
$ws.Range("E46").Value = "  -1.26%  "
$ws.Range("E47").Value = "  +6.65%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "26.54"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.34%  "
$ws.Range("E49").Value = "  +0.14%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "35.86"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +9.08%  "
$ws.Range("E51").Value = "  -0.58%  "
